# Insert a new data row at row 224 (pushes existing rows 224-340 down to 225-341)
# and populate it with a new price-report record for "Ají" (Americana (o) / Primera).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(224).Insert()

$ws.Cells.Item(224, 1).Value = 2
$ws.Cells.Item(224, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(224, 3).Value = "Coquimbo"
$ws.Cells.Item(224, 4).Value = 44846
$ws.Cells.Item(224, 5).Value = 4
$ws.Cells.Item(224, 6).Value = 100112021
$ws.Cells.Item(224, 7).Value = "Ají"
$ws.Cells.Item(224, 8).Value = "Americana (o)"
$ws.Cells.Item(224, 9).Value = "Primera"
$ws.Cells.Item(224, 10).Value = 300
$ws.Cells.Item(224, 11).Value = 45000
$ws.Cells.Item(224, 12).Value = 50000
$ws.Cells.Item(224, 13).Value = 47500
$ws.Cells.Item(224, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(224, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(224, 16).Value = 1900
$ws.Cells.Item(224, 17).Value = 25
$ws.Cells.Item(224, 18).Value = "Hortaliza"
